# Updated cryptos list (mirrors the GitHub Actions data refresh commit).
# Price (D) / Volume(1h) (E) values are refreshed for most rows; rows 44-45
# (Aave / MXToken) also swap back to their original name+link order with
# new price/volume figures.
#
# Some "Price" values (column D) are plain decimals that Excel would
# otherwise auto-parse as numbers (e.g. "217.43"); those are entered with a
# leading apostrophe to force text, matching the source data's string type,
# then the cell Style is reset to "Normal" so the apostrophe's quote-prefix
# formatting flag doesn't linger once the text value is committed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.917.23"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.638.60"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'217.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "'0.0623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'19.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.868.00"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "1.629.12"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "'67.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "26.900.83"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").Value = "  +2.37%  "
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").Value = "'9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'147.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "'7.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "'2.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "1.262.45"
$ws.Range("D36").Value = "'2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").Value = "'0.840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "'0.533"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'0.806"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "1.779.05"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'62.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("D46").Value = "'92.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  +18.18%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "'7.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
